# Update "Training Dashboard" sheet with new progress as of 04-Nov-2025:
#   - Column H (PERIOD TO EXPIRE) decreases by 1 day for every data row
#   - Column I (LAST UPDATE) changes from 03-Nov-2025 to 04-Nov-2025

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 37; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # Column H
    $iCell = $ws.Cells.Item($row, 9)   # Column I

    if ($iCell.Value2 -eq "03-Nov-2025") {
        $hCell.Value2 = $hCell.Value2 - 1
        $iCell.NumberFormat = "@"
        $iCell.Value2 = "04-Nov-2025"
    }
}
